{"js": "// Word JS API (Office.js) script.\n// Splits three \"run-on\" paragraphs (Programa PT, Programa EN, Bibliografia)\n// into multiple <w:t> segments joined by <w:br/> line breaks, keeping every\n// segment inside a single run (so run-level formatting, e.g. the italic\n// English paragraph, is preserved) \u2014 matching the target OOXML diff.\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build a <w:r> (optionally italic) containing alternating <w:t>/<w:br/>\n// children for the given list of text segments.\nfunction buildRunOoxml(segments, italic) {\n  const rPr = italic ? \"<w:rPr><w:i/></w:rPr>\" : \"\";\n  const pieces = segments.map((seg, i) => {\n    const preserve = /^\\s|\\s$/.test(seg) ? ' xml:space=\"preserve\"' : \"\";\n    const t = `<w:t${preserve}>${xmlEscape(seg)}</w:t>`;\n    return i < segments.length - 1 ? t + \"<w:br/>\" : t;\n  });\n  return `<w:r>${rPr}${pieces.join(\"\")}</w:r>`;\n}\n\nfunction wrapPackage(paragraphXml) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>${paragraphXml}</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n}\n\n// Split \"1. Foo.2. Bar.3. Baz.\" style run-on enumerations into separate\n// \"N. ...\" items by looking ahead for the next \"<digit>. \" marker.\nfunction splitNumberedList(text) {\n  return text.split(/(?=[0-9]+\\. )/).filter((s) => s.length > 0);\n}\n\n// The Bibliografia paragraph is a run-on list of distinct reference\n// entries; the break points are not algorithmically derivable from the\n// text alone (several entries start with \"SURNAME, \" for a co-author,\n// not a new entry), so the exact target entries are listed explicitly.\nconst BIBLIOGRAFIA_ENTRIES = [\n  \"CHIAVENATO, I. Introdu\u00e7\u00e3o \u00c0 Teoria Geral da Administra\u00e7\u00e3o. 9 ed. S\u00e3o Paulo: Manole, 2014.\",\n  \"CHIAVENATO, I. Administra\u00e7\u00e3o Para N\u00e3o Administradores: a Gest\u00e3o de Neg\u00f3cios Ao Alcance de Todos. 2 ed. S\u00e3o Paulo: Manole, 2011. \",\n  \"CHIAVENATO, I; SAPIRO, A. Planejamento Estrat\u00e9gico. Rio de Janeiro. Campus, 2004 \",\n  \"COLLINS, J.C.; PORRAS, J. I. Feitas para Durar: Pr\u00e1ticas bem-sucedidas de empresas vision\u00e1rias. 9\u00aa Ed.  Rio de Janeiro. Rocco, 2007 \",\n  \"GUERRINI, F. M.; ESCRI\u00c7\u00c3O FILHO, E.; ROSIM, D. Administra\u00e7\u00e3o Para Engenheiros. Rio de Janeiro: Campus, 2016.\",\n  \"HERRERO, E. Balanced Scorecard e a Gest\u00e3o Estrat\u00e9gica. Rio de Janeiro. Campus, 2005. \",\n  \"KAPLAN, R; NORTON, D. Kaplan e Norton na Pr\u00e1tica. Rio de Janeiro. Campus, 2004 \",\n  \"KAPLAN, R; NORTON, D. A Estrat\u00e9gia em A\u00e7\u00e3o: Balanced Scorecard. Rio de Janeiro. Campus, 1997 \",\n  \"KAPLAN, R; NORTON, D. Mapas Estrat\u00e9gicos. Rio de Janeiro. Campus, 2004 \",\n  \"MAXIMIANO, A. C. A. Teoria Geral da Administra\u00e7\u00e3o: da Revolu\u00e7\u00e3o Urbana \u00c0 Revolu\u00e7\u00e3o Digital. 8 ed. S\u00e3o Paulo: Atlas, 2017.\",\n  \"MINTZBERG, Henry; QUINN, James B. O processo da estrat\u00e9gia. 3\u00aa.ed.  Porto Alegre: Bookman, 2001.\",\n  \"MINTZBERG, H. Criando organiza\u00e7\u00f5es eficazes: estruturas em cinco configura\u00e7\u00f5es. 2\u00aa.ed. S\u00e3o Paulo: Atlas, 2003.\",\n  \"MORGAN, G. Imagens da organiza\u00e7\u00e3o. S\u00e3o Paulo, Atlas, 1996.\",\n  \"SILVA, M. M. L. Administra\u00e7\u00e3o para Estudantes e Profissionais de \u00c1reas T\u00e9cnicas. S\u00e3o Paulo: Brasport, 2018.\",\n  \"TZU, S. A Arte da Guerra (Edi\u00e7\u00e3o Completa). S\u00e3o Paulo. WMF Martins Fontes, 2009.\",\n];\n\nasync function run() {\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  let ptIndex = -1;\n  let enIndex = -1;\n  let biblioIndex = -1;\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const t = paragraphs.items[i].text;\n    if (ptIndex === -1 && t.indexOf(\"1. Teoria Geral de Administra\") === 0) {\n      ptIndex = i;\n    } else if (enIndex === -1 && t.indexOf(\"1. General Management Theory\") === 0) {\n      enIndex = i;\n    } else if (biblioIndex === -1 && t.indexOf(\"CHIAVENATO, I. Introdu\") === 0) {\n      biblioIndex = i;\n    }\n  }\n\n  if (ptIndex === -1 || enIndex === -1 || biblioIndex === -1) {\n    throw new Error(\n      \"Could not locate one or more target paragraphs (pt=\" +\n        ptIndex +\n        \" en=\" +\n        enIndex +\n        \" biblio=\" +\n        biblioIndex +\n        \")\"\n    );\n  }\n\n  // 1) Programa (Portuguese) \u2014 plain run.\n  const ptParagraph = paragraphs.items[ptIndex];\n  const ptSegments = splitNumberedList(ptParagraph.text);\n  const ptOoxml = wrapPackage(`<w:p>${buildRunOoxml(ptSegments, false)}</w:p>`);\n  ptParagraph.getRange().insertOoxml(ptOoxml, Word.InsertLocation.replace);\n  await context.sync();\n\n  // 2) Programa (English) \u2014 italic run.\n  const enParagraph = paragraphs.items[enIndex];\n  const enSegments = splitNumberedList(enParagraph.text);\n  const enOoxml = wrapPackage(`<w:p>${buildRunOoxml(enSegments, true)}</w:p>`);\n  enParagraph.getRange().insertOoxml(enOoxml, Word.InsertLocation.replace);\n  await context.sync();\n\n  // 3) Bibliografia \u2014 each entry separated by a blank line (two <w:br/>).\n  const biblioParagraph = paragraphs.items[biblioIndex];\n  const biblioPieces = BIBLIOGRAFIA_ENTRIES.map((seg, i) => {\n    const preserve = /^\\s|\\s$/.test(seg) ? ' xml:space=\"preserve\"' : \"\";\n    const t = `<w:t${preserve}>${xmlEscape(seg)}</w:t>`;\n    return i < BIBLIOGRAFIA_ENTRIES.length - 1 ? t + \"<w:br/><w:br/>\" : t;\n  }).join(\"\");\n  const biblioOoxml = wrapPackage(`<w:p><w:r>${biblioPieces}</w:r></w:p>`);\n  biblioParagraph.getRange().insertOoxml(biblioOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait run();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Splits three \"run-on\" paragraphs (Programa PT, Programa EN, Bibliografia)\n# into multiple <w:t> segments joined by <w:br/> line breaks, keeping every\n# segment inside a single run (so run-level formatting, e.g. the italic\n# English paragraph, is preserved) \u2014 matching the target OOXML diff.\n# Uses Range.InsertXML(...) so the exact <w:t>/<w:br/> shape (including\n# per-segment xml:space=\"preserve\") can be produced, rather than relying on\n# Range.Text assignment (which cannot reproduce per-segment xml:space).\n\nfunction XmlEscape($s) {\n    $s = $s -replace '&', '&amp;'\n    $s = $s -replace '<', '&lt;'\n    $s = $s -replace '>', '&gt;'\n    $s = $s -replace '\"', '&quot;'\n    return $s\n}\n\n# Split \"1. Foo.2. Bar.3. Baz.\" style run-on enumerations into separate\n# \"N. ...\" items by locating each \"<digit>. \" marker (regex here has no\n# look-around support, so boundaries are found via Matches + slicing).\nfunction Split-NumberedList($text) {\n    $ms = [regex]::Matches($text, '[0-9]+\\. ')\n    $parts = @()\n    $lastIndex = 0\n    $first = $true\n    foreach ($m in $ms) {\n        if ($first) {\n            # Keep the very first marker attached to the text that follows\n            # (do not split before index 0).\n            $first = $false\n            continue\n        }\n        $parts += $text.Substring($lastIndex, $m.Index - $lastIndex)\n        $lastIndex = $m.Index\n    }\n    $parts += $text.Substring($lastIndex)\n    return $parts\n}\n\n# Build a <w:r> (optionally italic) containing alternating <w:t>/<w:br/>\n# children for the given list of text segments, each gaining\n# xml:space=\"preserve\" only when it actually has leading/trailing\n# whitespace (matches the target diff exactly).\nfunction Build-RunOoxml($segments, $italic, $breaksBetween) {\n    $rpr = \"\"\n    if ($italic) {\n        $rpr = \"<w:rPr><w:i/></w:rPr>\"\n    }\n    $piecesList = @()\n    for ($i = 0; $i -lt $segments.Count; $i++) {\n        $seg = $segments[$i]\n        $preserve = \"\"\n        if ($seg -match '^\\s' -or $seg -match '\\s$') {\n            $preserve = ' xml:space=\"preserve\"'\n        }\n        $t = \"<w:t\" + $preserve + \">\" + (XmlEscape $seg) + \"</w:t>\"\n        if ($i -lt ($segments.Count - 1)) {\n            $t = $t + $breaksBetween\n        }\n        $piecesList += $t\n    }\n    $pieces = [string]::Join(\"\", $piecesList)\n    return \"<w:r>\" + $rpr + $pieces + \"</w:r>\"\n}\n\n# Range covering a paragraph's own text but NOT its trailing paragraph\n# mark. Using Paragraph.Range directly to drive InsertXML works for most\n# paragraphs, but when the paragraph is the very last one in the document\n# body, its Range end coincides with the end of the body and InsertXML\n# leaves a stray empty trailing paragraph behind; excluding the final\n# paragraph-mark character from the range avoids that in all cases.\nfunction Get-ParagraphContentRange($doc, $paragraph) {\n    $startPos = $paragraph.Range.Start\n    $endPos = $paragraph.Range.End - 1\n    return $doc.Range($startPos, $endPos)\n}\n\nfunction Wrap-Package($paragraphXml) {\n    $decl = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>'\n    $pkgOpen = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">'\n    $partOpen = '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>'\n    $docOpen = '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n    $docClose = '</w:body></w:document>'\n    $partClose = '</pkg:xmlData></pkg:part>'\n    $pkgClose = '</pkg:package>'\n    return $decl + $pkgOpen + $partOpen + $docOpen + $paragraphXml + $docClose + $partClose + $pkgClose\n}\n\n# The Bibliografia paragraph is a run-on list of distinct reference\n# entries; the break points are not algorithmically derivable from the\n# text alone (several entries start with \"SURNAME, \" for a co-author, not\n# a new entry), so the exact target entries are listed explicitly.\n$BibliografiaEntries = @(\n    \"CHIAVENATO, I. Introdu\u00e7\u00e3o \u00c0 Teoria Geral da Administra\u00e7\u00e3o. 9 ed. S\u00e3o Paulo: Manole, 2014.\",\n    \"CHIAVENATO, I. Administra\u00e7\u00e3o Para N\u00e3o Administradores: a Gest\u00e3o de Neg\u00f3cios Ao Alcance de Todos. 2 ed. S\u00e3o Paulo: Manole, 2011. \",\n    \"CHIAVENATO, I; SAPIRO, A. Planejamento Estrat\u00e9gico. Rio de Janeiro. Campus, 2004 \",\n    \"COLLINS, J.C.; PORRAS, J. I. Feitas para Durar: Pr\u00e1ticas bem-sucedidas de empresas vision\u00e1rias. 9\u00aa Ed.  Rio de Janeiro. Rocco, 2007 \",\n    \"GUERRINI, F. M.; ESCRI\u00c7\u00c3O FILHO, E.; ROSIM, D. Administra\u00e7\u00e3o Para Engenheiros. Rio de Janeiro: Campus, 2016.\",\n    \"HERRERO, E. Balanced Scorecard e a Gest\u00e3o Estrat\u00e9gica. Rio de Janeiro. Campus, 2005. \",\n    \"KAPLAN, R; NORTON, D. Kaplan e Norton na Pr\u00e1tica. Rio de Janeiro. Campus, 2004 \",\n    \"KAPLAN, R; NORTON, D. A Estrat\u00e9gia em A\u00e7\u00e3o: Balanced Scorecard. Rio de Janeiro. Campus, 1997 \",\n    \"KAPLAN, R; NORTON, D. Mapas Estrat\u00e9gicos. Rio de Janeiro. Campus, 2004 \",\n    \"MAXIMIANO, A. C. A. Teoria Geral da Administra\u00e7\u00e3o: da Revolu\u00e7\u00e3o Urbana \u00c0 Revolu\u00e7\u00e3o Digital. 8 ed. S\u00e3o Paulo: Atlas, 2017.\",\n    \"MINTZBERG, Henry; QUINN, James B. O processo da estrat\u00e9gia. 3\u00aa.ed.  Porto Alegre: Bookman, 2001.\",\n    \"MINTZBERG, H. Criando organiza\u00e7\u00f5es eficazes: estruturas em cinco configura\u00e7\u00f5es. 2\u00aa.ed. S\u00e3o Paulo: Atlas, 2003.\",\n    \"MORGAN, G. Imagens da organiza\u00e7\u00e3o. S\u00e3o Paulo, Atlas, 1996.\",\n    \"SILVA, M. M. L. Administra\u00e7\u00e3o para Estudantes e Profissionais de \u00c1reas T\u00e9cnicas. S\u00e3o Paulo: Brasport, 2018.\",\n    \"TZU, S. A Arte da Guerra (Edi\u00e7\u00e3o Completa). S\u00e3o Paulo. WMF Martins Fontes, 2009.\"\n)\n\n$d = $word.ActiveDocument\n\n$ptIndex = -1\n$enIndex = -1\n$biblioIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($ptIndex -eq -1 -and $t.StartsWith(\"1. Teoria Geral de Administra\")) {\n        $ptIndex = $i\n    } elseif ($enIndex -eq -1 -and $t.StartsWith(\"1. General Management Theory\")) {\n        $enIndex = $i\n    } elseif ($biblioIndex -eq -1 -and $t.StartsWith(\"CHIAVENATO, I. Introdu\")) {\n        $biblioIndex = $i\n    }\n    $i++\n}\n\nif ($ptIndex -eq -1 -or $enIndex -eq -1 -or $biblioIndex -eq -1) {\n    throw \"Could not locate one or more target paragraphs (pt=$ptIndex en=$enIndex biblio=$biblioIndex)\"\n}\n\n# 1) Programa (Portuguese) \u2014 plain run, single <w:br/> between items.\n$ptParagraph = $d.Paragraphs($ptIndex)\n$ptText = $ptParagraph.Range.Text\n$ptText = $ptText.TrimEnd([char]13, [char]7)\n$ptSegments = Split-NumberedList $ptText\n$ptRun = Build-RunOoxml $ptSegments $false \"<w:br/>\"\n$ptRange = Get-ParagraphContentRange $d $ptParagraph\n$null = $ptRange.InsertXML((Wrap-Package (\"<w:p>\" + $ptRun + \"</w:p>\")))\n\n# 2) Programa (English) \u2014 italic run, single <w:br/> between items.\n$enParagraph = $d.Paragraphs($enIndex)\n$enText = $enParagraph.Range.Text\n$enText = $enText.TrimEnd([char]13, [char]7)\n$enSegments = Split-NumberedList $enText\n$enRun = Build-RunOoxml $enSegments $true \"<w:br/>\"\n$enRange = Get-ParagraphContentRange $d $enParagraph\n$null = $enRange.InsertXML((Wrap-Package (\"<w:p>\" + $enRun + \"</w:p>\")))\n\n# 3) Bibliografia \u2014 each entry separated by a blank line (two <w:br/>).\n$biblioParagraph = $d.Paragraphs($biblioIndex)\n$biblioRun = Build-RunOoxml $BibliografiaEntries $false \"<w:br/><w:br/>\"\n$biblioRange = Get-ParagraphContentRange $d $biblioParagraph\n$null = $biblioRange.InsertXML((Wrap-Package (\"<w:p>\" + $biblioRun + \"</w:p>\")))\n"}
